$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data rows (2-5) each have their Fecha/Volumen/Precio fields
# cyclically shifted: row N now holds what used to be in row N+1 (row 5
# wraps around to what used to be in row 2) for columns D, J, K, L, M, P.

# Capture the "before" values first so the shift can be applied atomically.
$rows = 2..5
$before = @{}
foreach ($r in $rows) {
    $before[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        P = $ws.Cells.Item($r, 16).Value2
    }
}

$map = @{ 2 = 3; 3 = 4; 4 = 5; 5 = 2 }

foreach ($r in $rows) {
    $src = $map[$r]
    $vals = $before[$src]
    $ws.Cells.Item($r, 4).Value2 = $vals.D
    $ws.Cells.Item($r, 10).Value2 = $vals.J
    $ws.Cells.Item($r, 11).Value2 = $vals.K
    $ws.Cells.Item($r, 12).Value2 = $vals.L
    $ws.Cells.Item($r, 13).Value2 = $vals.M
    $ws.Cells.Item($r, 16).Value2 = $vals.P
}

$wb.Save()
